$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "Weekdays Total" (column B) and "Weekends Total" (column C) values
$ws.Range("B2").Value = 10
$ws.Range("C2").Value = 7

$ws.Range("B3").Value = 9
$ws.Range("C3").Value = 7

$ws.Range("B4").Value = 9
$ws.Range("C4").Value = 7

$ws.Range("B5").Value = 9
$ws.Range("C5").Value = 6

$ws.Range("B6").Value = 9
$ws.Range("C6").Value = 8

$ws.Range("B7").Value = 9
$ws.Range("C7").Value = 7

$ws.Range("B8").Value = 9
$ws.Range("C8").Value = 6

$ws.Range("B9").Value = 9
$ws.Range("C9").Value = 7

$ws.Range("B10").Value = 9
$ws.Range("C10").Value = 7

$ws.Range("B11").Value = 9
$ws.Range("C11").Value = 7

$ws.Range("B12").Value = 9
$ws.Range("C12").Value = 7

$ws.Range("B13").Value = 9
$ws.Range("C13").Value = 7

$ws.Range("B14").Value = 9
$ws.Range("C14").Value = 8

$ws.Range("B15").Value = 9
$ws.Range("C15").Value = 7

$ws.Range("B16").Value = 9
$ws.Range("C16").Value = 7
